# Updates cryptos list price/volume columns (and the Hedera/WEMIXToken
# row swap at rows 32-33) to match the latest scrape, per the GitHub Actions
# "Updated cryptos list" commit.
#
# Price/volume cells are plain text in this sheet (e.g. "35.171.57",
# "  -0.36%  "), not numbers. Excel's COM Range.Value setter auto-detects
# numeric-looking strings (e.g. "253.08") and coerces them to real numbers,
# which would corrupt values like "18.50" (-> 18.5) or introduce float noise
# (-> 253.08000000000001). A leading apostrophe is Excel's standard "force
# text" input prefix, so it's used for the ambiguous numeric-looking D-column
# values to keep them exact text, matching the target content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.172.88"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.900.89"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'253.08"
$ws.Range("E5").Value = "  +2.91%  "
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'40.96"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").Value = "'0.359"
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("D10").Value = "'52.85"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").Value = "'0.0753"
$ws.Range("E11").Value = "  +3.77%  "
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("E13").Value = "  +4.73%  "
$ws.Range("D14").Value = "2.176.08"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "'0.735"
$ws.Range("E15").Value = "  +3.23%  "
$ws.Range("D16").Value = "'4.97"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("D17").Value = "1.900.06"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "35.148.89"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "'73.61"
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").Value = "0.0₃0834"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").Value = "'242.59"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "'12.98"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("E23").Value = "  +4.90%  "
$ws.Range("E25").Value = "  +4.56%  "
$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'167.19"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").Value = "'8.58"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'18.50"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").Value = "4.128.49"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.0614"
$ws.Range("E32").Value = "  +7.45%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.04"
$ws.Range("E33").Value = "  +14.64%  "
$ws.Range("D34").Value = "'4.31"
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("E35").Value = "  +7.80%  "
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'0.852"
$ws.Range("E38").Value = "  -11.23%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").Value = "'102.52"
$ws.Range("E40").Value = "  +13.90%  "
$ws.Range("D41").Value = "'17.41"
$ws.Range("E41").Value = "  +8.03%  "
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'0.0646"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("D45").Value = "'2.43"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").Value = "1.321.86"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").Value = "'6.59"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").Value = "'11.94"
$ws.Range("E50").Value = "  -6.32%  "
$ws.Range("D51").Value = "'43.11"
$ws.Range("E51").Value = "  -6.64%  "
